# Indexant llibre Blancafort anys 1919 - 1952
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column P ("Estat"), pushing it to column Q.
$colP = $ws.Columns.Item(16)
$colP.Insert()

# Give the new "Llibre" column (P) the same visual width as its neighbour "Matrimoni Pares" (O).
$ws.Columns.Item(16).ColumnWidth = 18.417

# New column P header + values: "Llibre" = book/volume reference for every data row.
$ws.Cells.Item(1, 16).Value = "Llibre"
$bookName = "6 Baptismes 1919-1952"
for ($r = 2; $r -le 47; $r++) {
    $ws.Cells.Item($r, 16).Value = $bookName
}

# Restore the view (scroll position / selection) to match the saved state.
$ws.Range("A6").Select()
